# ============================================================
# Trade #109 closed at 2026-02-16 21:42:45 - leadlag DOWN +0.000%
#
# This applies:
#  1) Summary sheet stat refresh (overall + leadlag totals)
#  2) leadlag sheet: trade #85 (row 65) marked CLOSED with exit data,
#     and a brand-new open trade #109 appended as row 85
#  3) All Trades sheet: mirrored CLOSED record appended as row 86
#  4) Comparison sheet: leadlag aggregate stat refresh
# ============================================================

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    # Forces a cell to literal text (so date-like / percent-like / numeric-
    # looking strings such as "2026-02-16", "69.4%" or "2.83" are stored
    # verbatim instead of being auto-coerced into a number), then restores
    # the cell's style to the workbook default so only the content differs.
    param($Range, $Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

# ============================================================
# Sheet: Summary
# ============================================================
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("C2").Value = 85
Set-TextCell $wsSummary.Range("D2") "69.4%"
Set-TextCell $wsSummary.Range("E2") "+24.7055%"
Set-TextCell $wsSummary.Range("F2") "+0.2907%"

$wsSummary.Range("C3").Value = 83
Set-TextCell $wsSummary.Range("D3") "48.2%"
Set-TextCell $wsSummary.Range("E3") "+13.3842%"
Set-TextCell $wsSummary.Range("F3") "+0.1613%"

Write-Host "Summary sheet updated"

# ============================================================
# Sheet: leadlag
# ============================================================
$wsLead = $wb.Worksheets.Item("leadlag")

# --- Existing trade (row 65) transitions from OPEN to CLOSED ---
$wsLead.Range("G65").Value = 68561.25019799999
$wsLead.Range("H65").Value = "CLOSED"
$wsLead.Range("I65").Value = -0.0414
$wsLead.Range("J65").Value = -0.41
Set-TextCell $wsLead.Range("M65") "time_exit_5min"
$wsLead.Range("N65").Value = 5

# --- New trade appended as row 85 (trade #109, still OPEN) ---
$wsLead.Range("A85").Value = 109
Set-TextCell $wsLead.Range("B85") "2026-02-16"
Set-TextCell $wsLead.Range("C85") "21:42:45"
Set-TextCell $wsLead.Range("D85") "leadlag"
Set-TextCell $wsLead.Range("E85") "DOWN"
$wsLead.Range("F85").Value = 68395.645
Set-TextCell $wsLead.Range("H85") "OPEN"
$wsLead.Range("I85").Value = 0
$wsLead.Range("J85").Value = 0
$wsLead.Range("K85").Value = 0.75
Set-TextCell $wsLead.Range("L85") "Binance leading with -0.107% move"
$wsLead.Range("N85").Value = 0

Write-Host "leadlag sheet updated"

# ============================================================
# Sheet: All Trades
# ============================================================
$wsAll = $wb.Worksheets.Item("All Trades")

# --- Mirror of the now-CLOSED leadlag trade, appended as row 86 ---
$wsAll.Range("A86").Value = 85
Set-TextCell $wsAll.Range("B86") "2026-02-16"
Set-TextCell $wsAll.Range("C86") "21:37:42"
Set-TextCell $wsAll.Range("D86") "leadlag"
Set-TextCell $wsAll.Range("E86") "DOWN"
$wsAll.Range("F86").Value = 68532.86500000001
$wsAll.Range("G86").Value = 68561.25019799999
Set-TextCell $wsAll.Range("H86") "CLOSED"
$wsAll.Range("I86").Value = -0.0414
$wsAll.Range("J86").Value = -0.41
$wsAll.Range("K86").Value = 0.75
Set-TextCell $wsAll.Range("L86") "Binance leading with -0.132% move"
Set-TextCell $wsAll.Range("M86") "time_exit_5min"
$wsAll.Range("N86").Value = 5

Write-Host "All Trades sheet updated"

# ============================================================
# Sheet: Comparison
# ============================================================
$wsComp = $wb.Worksheets.Item("Comparison")

$wsComp.Range("B2").Value = 83
Set-TextCell $wsComp.Range("C2") "48.2%"
Set-TextCell $wsComp.Range("D2") "2.83"
Set-TextCell $wsComp.Range("F2") "-0.3052%"
Set-TextCell $wsComp.Range("G2") "1.70"

Write-Host "Comparison sheet updated"

Write-Host "All updates complete"
